$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wynagrodzenie")

# --- Step 1: freeze the old helper formulas (columns D "zmiana %" and E "okres") to plain values
# so that inserting/deleting/sorting columns below doesn't corrupt their relative references.
$ws.Range("D2:E71").Value = $ws.Range("D2:E71").Value()

# --- Step 2: the old column E ("YYYY-12" period text) is being replaced by a brand-new column B
# with the same kind of text, so drop the old E column now.
$ws.Columns.Item(5).Delete()

# --- Step 3: insert a fresh, blank column B. This shifts:
#       old B (category)      -> C
#       old C (value)         -> D
$ws.Columns.Item(2).Insert()

# --- Step 4: header row
$ws.Range("B1").Value = "dr"

# --- Step 5: fill the new column B with the "YYYY-12" period label for every data row
for ($r = 2; $r -le 71; $r++) {
    $year = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = "$year-12"
}

# --- Step 6: re-sort the whole data block (rows 2-71) by category (now column C), descending -
# matches the original table's sort order (sortCondition descending on the category column).
$sortRange = $ws.Range("A2:D71")
$sortKey = $ws.Range("C2:C71")
$sortRange.Sort($sortKey, 2, $null, $null, 1, $null, $null, 2)

# --- Step 7: rebuild the "zmiana %" formulas in column E for every row except the last row of
# each 10-row (2012-2021) block, referencing the now-correct column D.
for ($base = 2; $base -le 62; $base = $base + 10) {
    for ($i = 0; $i -le 8; $i++) {
        $r = $base + $i
        $rNext = $r + 1
        $ws.Cells.Item($r, 5).Formula = "=(D$r-D$rNext)/D$rNext*100"
    }
}

# --- Step 8: rebuild the AutoFilter so it only spans the header row, and update the matching
# workbook-level _FilterDatabase defined name to follow.
$ws.AutoFilterMode = $false
$ws.Range("A1:E1").AutoFilter()
$filterName = $wb.Names.Item("Wynagrodzenie!_FilterDatabase")
$filterName.RefersTo = "=Wynagrodzenie!`$A`$1:`$E`$1"

# --- Step 9: move the selection / active cell on this sheet, and switch the active tab to
# "Inflacja1" (activating it moves Excel's tabSelected flag there and updates workbook activeTab).
$ws.Range("K16:K17").Select()
$wsInflacja1 = $wb.Worksheets.Item("Inflacja1")
$wsInflacja1.Activate()
